$d = $word.ActiveDocument

# Build the OOXML for the two new paragraphs that need to be inserted at the
# very top of the document body:
#   1) a paragraph with two runs (bio blurb)
#   2) a completely empty paragraph (spacer)
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$run1Text = "I" + [char]0x2019 + "m Richard Fu, an artist, designer, and developer who loves blue"
$run2Text = ". This is my portfolio website. It contains some of my most recent projects as well as ways to reach me. "

$xml = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t>' + $run1Text + '</w:t></w:r><w:r><w:t xml:space="preserve">' + $run2Text + '</w:t></w:r></w:p><w:p xmlns:w="' + $wNs + '"/>'

$insertPoint = $d.Range(0, 0)
$insertPoint.InsertXML($xml)
